$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.400.71"
$ws.Range("E2").Value = "  -0.56%  "

$ws.Range("D3").Value = "1.829.22"
$ws.Range("E3").Value = "  +1.90%  "

$ws.Range("E4").Value = "  -0.30%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "330.17"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.26%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4477"
$ws.Range("E7").Value = "  +1.93%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3788"
$ws.Range("E8").Value = "  +1.60%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "44.71"
$ws.Range("E9").Value = "  -2.21%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07787"
$ws.Range("E10").Value = "  +2.86%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.143"
$ws.Range("E11").Value = "  +1.07%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.49"
$ws.Range("E12").Value = "  -0.29%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.002"
$ws.Range("E13").Value = "  -0.27%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.381"
$ws.Range("E14").Value = "  +3.02%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.543"
$ws.Range("E15").Value = "  +0.99%  "

$ws.Range("D16").Value = "1.838.40"
$ws.Range("E16").Value = "  +2.16%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "93.52"
$ws.Range("E17").Value = "  +16.12%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001088"
$ws.Range("E18").Value = "  +0.31%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06411"
$ws.Range("E19").Value = "  -4.37%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.001"
$ws.Range("E20").Value = "  -0.19%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.59"
$ws.Range("E21").Value = "  +0.73%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.376"
$ws.Range("E22").Value = "  +2.57%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.5428"
$ws.Range("E23").Value = "  -0.26%  "

$ws.Range("D24").Value = "28.462.19"
$ws.Range("E24").Value = "  -0.41%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.77"
$ws.Range("E25").Value = "  +0.96%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.270"
$ws.Range("E26").Value = "  -6.73%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.88"
$ws.Range("E27").Value = "  +2.31%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "154.43"
$ws.Range("E28").Value = "  +0.81%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.383"
$ws.Range("E29").Value = "  +2.23%  "

$ws.Range("D30").Value = "2.046.08"
$ws.Range("E30").Value = "  +1.95%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "129.09"
$ws.Range("E31").Value = "  -1.05%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.217"
$ws.Range("E32").Value = "  -7.13%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.942"
$ws.Range("E33").Value = "  +3.06%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.09307"
$ws.Range("E34").Value = "  +1.08%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.671"
$ws.Range("E35").Value = "  -7.71%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "13.13"
$ws.Range("E36").Value = "  +8.90%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02366"
$ws.Range("E37").Value = "  +2.63%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2207"
$ws.Range("E38").Value = "  -1.22%  "

$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6689"
$ws.Range("E39").Value = "  +1.89%  "

$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06315"
$ws.Range("E40").Value = "  +0.84%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.232"
$ws.Range("E41").Value = "  +1.28%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.199"
$ws.Range("E42").Value = "  +2.59%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.199"
$ws.Range("E43").Value = "  +0.34%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.001"
$ws.Range("E44").Value = "  -0.13%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.09"
$ws.Range("E45").Value = "  +2.56%  "

$ws.Range("B46").Value = "WEMIXTOKEN"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.405"
$ws.Range("E46").Value = "  -1.48%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6177"
$ws.Range("E47").Value = "  +2.11%  "

$ws.Range("E48").Value = "  -0.28%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.066"
$ws.Range("E49").Value = "  +3.01%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "127.72"
$ws.Range("E50").Value = "  +0.38%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07015"
$ws.Range("E51").Value = "  +0.20%  "
